$wb = $excel.ActiveWorkbook

$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"

# "Status" column on each language sheet shares the same string as the
# Overview sheet's status cells ("Ready for handoff" -> "Handback transform failed")
$zhSheet.Range("C3").Value = "Handback transform failed"
$deSheet.Range("C3").Value = "Handback transform failed"

$zhSheet.Range("K3").Value = "Handback file name: natzk5wk.nnc is different with handoff file name: 7c390d2b-e229-477a-a819-d1abccae1ca4.1957febbefdc96a539a848e535dd56ddcbb8bedb.zh-cn."

$deSheet.Range("K3").Value = "Handback file name: natzk5wk.nnc is different with handoff file name: 7c390d2b-e229-477a-a819-d1abccae1ca4.1957febbefdc96a539a848e535dd56ddcbb8bedb.de-de."
